$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planned BOM")

# The old row 3 ("Sparkfun Qwiic Pro Micro" microcontroller entry) is removed entirely.
# Deleting it shifts the old row 4 ("Stranded Wire" / "Wire Limit Switch") up to row 3,
# which is exactly the row 3 content of the updated BOM.
$ws.Rows.Item(3).Delete()

# New row 4: Cam Follower / Notch Detents, with a hyperlinked Amazon link.
$ws.Range("D4").Value = "https://www.amazon.com/gp/product/B08C5CTNZR/"
$ws.Range("A4").Value = "Cam Follower"
$ws.Range("B4").Value = "Notch Detents"
$ws.Range("C4").Value = 1

# New row 5: Ball Bearings Mounted / Bearing support on wall, plain-text Amazon link.
$ws.Range("D5").Value = "https://www.amazon.com/gp/product/B07C5NPXMZ/"
$ws.Range("A5").Value = "Ball Bearings Mounted"
$ws.Range("B5").Value = "Bearing support on wall"
$ws.Range("C5").Value = 1

# Turn the D4 link cell into a real hyperlink (adds the Hyperlink cell style too).
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.amazon.com/gp/product/B08C5CTNZR/")

# Column B/D got wider once the longer "Purpose"/"Link" text was auto-fit.
$ws.Columns.Item(2).ColumnWidth = 21.5
$ws.Columns.Item(4).ColumnWidth = 97

# Matches the cursor position recorded in the saved file.
$ws.Range("D13").Select() | Out-Null
